# Data-content edit extracted from the commit diff:
# The shared string used in cell A1 changes from "Grade" to "grade".
# (All other diff hunks are OOXML re-save metadata/noise produced by a
# newer Excel build re-serializing the package - fileVersion, revisionPtr,
# xmlns additions, sheetView selection, dyDescent rounding, and the
# customXml part renumbering - none of which represent a workbook content
# change, so they are not reproduced here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "grade"
